# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values for each data row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 0
    7  = 4
    8  = 0
    9  = 0
    10 = 4
    11 = 2
    12 = 3
    13 = 0
    14 = 2
    15 = 0
    16 = 1
    17 = 0
    18 = 0
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 2
    27 = 1
    28 = 1
    29 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
